$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E3").Value = "SHP"
$ws.Range("J3").Value = "e"
$ws.Range("K3").Value = "e"
$ws.Range("L3").Value = "e"

$ws.Range("D5").Value = "n"
$ws.Range("E5").Value = "FAC"
$ws.Range("F5").Value = "n"
